$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 424.2
$ws.Range("I18").Value = 514.25
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 514.25
$ws.Range("L18").Value = 64
$ws.Range("M18").Value = -230.25
$ws.Range("N18").Value = -632
$ws.Range("H74").Value = 5673.3335
$ws.Range("I74").Value = 5020
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 5020
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -4084
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 5673.3335
$ws.Range("I77").Value = 5020
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 25100
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -20420
$ws.Range("N77").Value = -39360
$ws.Range("H125").Value = 8509.625
$ws.Range("I125").Value = 11000
$ws.Range("J125").Value = 7679.5
$ws.Range("K125").Value = 99000
$ws.Range("L125").Value = 69115.5
$ws.Range("M125").Value = -96540
$ws.Range("N125").Value = -74035.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5259.143
$ws.Range("I2").Value = 3223.75
$ws.Range("J2").Value = 7973
$ws.Range("K2").Value = 3223.75
$ws.Range("L2").Value = 7973
$ws.Range("M2").Value = -3110.75
$ws.Range("N2").Value = -8199
$ws.Range("H22").Value = 2100
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1801
$ws.Range("H32").Value = 4166.643
$ws.Range("I32").Value = 3953.508
$ws.Range("K32").Value = 3953.508
$ws.Range("M32").Value = -3666.508
$ws.Range("H45").Value = 3250.5
$ws.Range("J45").Value = 4652
$ws.Range("L45").Value = 4652
$ws.Range("N45").Value = -5406
$ws.Range("H74").Value = 35023.055
$ws.Range("I74").Value = 38571.363
$ws.Range("J74").Value = 5749.5
$ws.Range("K74").Value = 38571.363
$ws.Range("L74").Value = 5749.5
$ws.Range("M74").Value = -37697.363
$ws.Range("N74").Value = -7497.5
$ws.Range("H77").Value = 35023.055
$ws.Range("I77").Value = 38571.363
$ws.Range("J77").Value = 5749.5
$ws.Range("K77").Value = 192856.815
$ws.Range("L77").Value = 28747.5
$ws.Range("M77").Value = -188488.815
$ws.Range("N77").Value = -37483.5
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H102").Value = 3476.6428
$ws.Range("I102").Value = 2993.4348
$ws.Range("K102").Value = 2993.4348
$ws.Range("M102").Value = -1371.4348
$ws.Range("H116").Value = 5259.143
$ws.Range("I116").Value = 3223.75
$ws.Range("J116").Value = 7973
$ws.Range("K116").Value = 3223.75
$ws.Range("L116").Value = 7973
$ws.Range("M116").Value = -929.75
$ws.Range("N116").Value = -12561

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5259.143
$ws.Range("I3").Value = 3223.75
$ws.Range("J3").Value = 7973
$ws.Range("K3").Value = 3223.75
$ws.Range("L3").Value = 7973
$ws.Range("M3").Value = -3109.75
$ws.Range("N3").Value = -8201
$ws.Range("H20").Value = 1900.3636
$ws.Range("J20").Value = 2354.75
$ws.Range("L20").Value = 2354.75
$ws.Range("N20").Value = -2848.75
$ws.Range("H22").Value = 400
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -746
$ws.Range("H58").Value = 40348
$ws.Range("I58").Value = 32111
$ws.Range("J58").Value = 48585
$ws.Range("K58").Value = 32111
$ws.Range("L58").Value = 48585
$ws.Range("M58").Value = -31817
$ws.Range("N58").Value = -49173
$ws.Range("H92").Value = 243827.33
$ws.Range("J92").Value = 243827.33
$ws.Range("L92").Value = 243827.33
$ws.Range("N92").Value = -248819.33
$ws.Range("H94").Value = 1516.7273
$ws.Range("I94").Value = 1000.4545
$ws.Range("K94").Value = 1000.4545
$ws.Range("M94").Value = -549.4545000000001
$ws.Range("H99").Value = 4527.778
$ws.Range("I99").Value = 2809.6667
$ws.Range("K99").Value = 2809.6667
$ws.Range("M99").Value = -1311.6667
$ws.Range("H107").Value = 2676.5667
$ws.Range("I107").Value = 2458.1428
$ws.Range("J107").Value = 5734.5
$ws.Range("K107").Value = 2458.1428
$ws.Range("L107").Value = 5734.5
$ws.Range("M107").Value = -538.1428000000001
$ws.Range("N107").Value = -9574.5
$ws.Range("H134").Value = 898.375
$ws.Range("I134").Value = 898.375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2695.125
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -160.125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2149.75
$ws.Range("I16").Value = 2349.5
$ws.Range("K16").Value = 2349.5
$ws.Range("M16").Value = -2062.5
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20472
$ws.Range("H30").Value = 20000
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20182
$ws.Range("H93").Value = 20973.615
$ws.Range("I93").Value = 18181
$ws.Range("J93").Value = 30282.334
$ws.Range("K93").Value = 18181
$ws.Range("L93").Value = 30282.334
$ws.Range("M93").Value = -16309
$ws.Range("N93").Value = -34026.334
$ws.Range("H113").Value = 2149.75
$ws.Range("I113").Value = 2349.5
$ws.Range("K113").Value = 2349.5
$ws.Range("M113").Value = -179.5
$ws.Range("H128").Value = 20000
$ws.Range("J128").Value = 20000
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 630.4
$ws.Range("I68").Value = 550.6667
$ws.Range("J68").Value = 750
$ws.Range("K68").Value = 1652.0001
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -841.0001
$ws.Range("N68").Value = -3872
$ws.Range("H71").Value = 630.4
$ws.Range("I71").Value = 550.6667
$ws.Range("J71").Value = 750
$ws.Range("K71").Value = 4956.0003
$ws.Range("L71").Value = 6750
$ws.Range("M71").Value = -900.0002999999997
$ws.Range("N71").Value = -14862
$ws.Range("H115").Value = 4133.4
$ws.Range("J115").Value = 6619.8335
$ws.Range("L115").Value = 19859.5005
$ws.Range("N115").Value = -22209.5005
$ws.Range("H140").Value = 1920.3529
$ws.Range("I140").Value = 1467.6428
$ws.Range("J140").Value = 4033
$ws.Range("K140").Value = 4402.928400000001
$ws.Range("L140").Value = 12099
$ws.Range("M140").Value = 777.0715999999993
$ws.Range("N140").Value = -22459

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 171.6
$ws.Range("I2").Value = 145.5
$ws.Range("K2").Value = 145.5
$ws.Range("M2").Value = -32.5
$ws.Range("H114").Value = 65833.336
$ws.Range("J114").Value = 65833.336
$ws.Range("L114").Value = 65833.336
$ws.Range("N114").Value = -74511.336
$ws.Range("H122").Value = 1364.6875
$ws.Range("I122").Value = 1089
$ws.Range("K122").Value = 3267
$ws.Range("M122").Value = -817
$ws.Range("H132").Value = 2343.7273
$ws.Range("I132").Value = 1917.1875
$ws.Range("J132").Value = 3481.1667
$ws.Range("K132").Value = 5751.5625
$ws.Range("L132").Value = 10443.5001
$ws.Range("M132").Value = -3221.5625
$ws.Range("N132").Value = -15503.5001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4476.5386
$ws.Range("I7").Value = 3578.2222
$ws.Range("K7").Value = 3578.2222
$ws.Range("M7").Value = -3466.2222
$ws.Range("H122").Value = 4029.75
$ws.Range("I122").Value = 3723.0833
$ws.Range("K122").Value = 11169.2499
$ws.Range("M122").Value = -8719.249899999999
$ws.Range("H126").Value = 4476.5386
$ws.Range("I126").Value = 3578.2222
$ws.Range("K126").Value = 10734.6666
$ws.Range("M126").Value = -8264.6666
$ws.Range("H132").Value = 6159.268
$ws.Range("I132").Value = 1421.4546
$ws.Range("J132").Value = 25702.75
$ws.Range("K132").Value = 4264.3638
$ws.Range("L132").Value = 77108.25
$ws.Range("M132").Value = -1734.3638
$ws.Range("N132").Value = -82168.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8112.375
$ws.Range("I62").Value = 6499
$ws.Range("K62").Value = 6499
$ws.Range("M62").Value = -5875
$ws.Range("H65").Value = 8112.375
$ws.Range("I65").Value = 6499
$ws.Range("K65").Value = 32495
$ws.Range("M65").Value = -29375
$ws.Range("H132").Value = 1009.5357
$ws.Range("I132").Value = 985.4091
$ws.Range("J132").Value = 1098
$ws.Range("K132").Value = 2956.2273
$ws.Range("L132").Value = 3294
$ws.Range("M132").Value = -426.2273
$ws.Range("N132").Value = -8354

Write-Output "applied changes"
